$wb = $excel.ActiveWorkbook
$wsLib = $wb.Worksheets.Item("libraries")
$wsCsv = $wb.Worksheets.Item("save to libraries.csv")

# Fix the description text: replace the unicode ellipsis character with three literal dots
$wsLib.Range("B3").Value = "exporting images to SVG, EPS, ..."

# Set explicit column widths on the "save to libraries.csv" sheet (A=30 chars, B=39 chars)
$wsCsv.Columns.Item(1).ColumnWidth = 29.166666666666668
$wsCsv.Columns.Item(2).ColumnWidth = 38.166666666666664

# Update selection/cursor state on both sheets, and make "libraries" the active tab
$wsCsv.Range("B2").Select() | Out-Null
$wsLib.Activate() | Out-Null
$wsLib.Range("B3").Select() | Out-Null
